$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Net Income
$ws.Range("B2").Value = 2761138000.0

# Row 3 - Depreciation & Amortization
$ws.Range("B3").Value = 111528000.0

# Row 4 - Non Cash Items (Other)
$ws.Range("B4").Value = 501192000.0

# Row 5 - Accounts Receivable Change
$ws.Range("B5").Value = -98145000.0

# Row 6 - Change in inventories
$ws.Range("B6").Value = -127349000.0
$ws.Range("C6").Value = -132014000.0
$ws.Range("D6").Value = -113674000.0
$ws.Range("E6").Value = -103259000.0
$ws.Range("F6").Value = -78358000.0
$ws.Range("G6").Value = -64047000.0

# Row 7 - Accounts Payable Change
$ws.Range("B7").Value = 26347000.0

# Row 8 - Change in payables and accrued liability
$ws.Range("B8").Value = 516983000.0
$ws.Range("C8").Value = 465000000.0
$ws.Range("D8").Value = 398000000.0
$ws.Range("E8").Value = 431731000.0
$ws.Range("F8").Value = 471978000.0
$ws.Range("G8").Value = 343699000.0

# Row 10 - Change in other assets and liabilities
$ws.Range("B10").Value = 61890000.0

# Row 11 - Operating Cash Flow
$ws.Range("B11").Value = 3358862000.0

# Row 12 - Capital expenditures
$ws.Range("B12").Value = -311274000.0

# Row 14 - Long-Term Investments Change (Net)
$ws.Range("B14").Value = 304687000.0

# Row 16 - Investing cash flow
$ws.Range("B16").Value = -6587000.0

# Row 17 - Repayment/Issuance of Debt (Net)
$ws.Range("B17").Value = -33171000.0

# Row 18 - Equity Repurchase (Common, Net)
$ws.Range("B18").Value = -463154000.0

# Row 19 - Other financial activities
$ws.Range("B19").Value = -130568000.0

# Row 20 - Financing cash flow
$ws.Range("B20").Value = -664505000.0

# Row 21 - Exchange Rate Adjustment
$ws.Range("B21").Value = 23173000.0

# Row 22 - Change in Cash
$ws.Range("B22").Value = 2710943000.0

# Row 23 - Beginning Cash
$ws.Range("B23").Value = 5988845000.0

# Row 24 - Ending Cash
$ws.Range("B24").Value = 8699788000.0

# Row 25 - Stock Based Compensation
$ws.Range("B25").Value = 428929000.0

# Row 27 - Assets Liabilities Change (Total)
$ws.Range("B27").Value = -15890000.0

# Row 28 - Investments Change (Net)
$ws.Range("B28").Value = 304687000.0

# Row 29 - Issuance/Purchase of Shares
$ws.Range("B29").Value = -463154000.0
